$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row (B11): correct/total marks adjustment
$ws.Range("B11").Value = 5

# Update the "Total" row (B12) and the corr/total display (E12)
$ws.Range("B12").Value = 105
$ws.Range("E12").Value = "105/140"
